# The workbook's single data sheet gained two brand-new observation rows
# (a weekly refresh of "Poroto verde" price data). In the published sheet
# these show up at rows 715-716, pushing every existing data row down by
# two positions (old row 715 -> new row 717, ... old row 812 -> new row 814),
# and the sheet's used range grows from A1:R812 to A1:R814.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top of where the new data belongs.
# Inserting the two rows together shifts every row at/after 715 down by 2,
# exactly matching the before/after row alignment implied by the diff.
$ws.Range("A715:A716").EntireRow.Insert()

# New row 715: Magnum / Primera, Peru origin, $/malla 25 kilos
$ws.Range("A715").Value = 9
$ws.Range("B715").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C715").Value = "Metropolitana"
$ws.Range("D715").Value = 45127
$ws.Range("E715").Value = 13
$ws.Range("F715").Value = 100112031
$ws.Range("G715").Value = "Poroto verde"
$ws.Range("H715").Value = "Magnum"
$ws.Range("I715").Value = "Primera"
$ws.Range("J715").Value = 70
$ws.Range("K715").Value = 19000
$ws.Range("L715").Value = 21000
$ws.Range("M715").Value = 20000
$ws.Range("N715").Value = "$/malla 25 kilos"
$ws.Range("O715").Value = "Perú"
$ws.Range("P715").Value = 800
$ws.Range("Q715").Value = 25
$ws.Range("R715").Value = "Hortaliza"

# New row 716: Sin especificar / Primera, Peru origin, $/malla 25 kilos
$ws.Range("A716").Value = 9
$ws.Range("B716").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C716").Value = "Metropolitana"
$ws.Range("D716").Value = 45127
$ws.Range("E716").Value = 13
$ws.Range("F716").Value = 100112031
$ws.Range("G716").Value = "Poroto verde"
$ws.Range("H716").Value = "Sin especificar"
$ws.Range("I716").Value = "Primera"
$ws.Range("J716").Value = 52
$ws.Range("K716").Value = 20000
$ws.Range("L716").Value = 22000
$ws.Range("M716").Value = 21000
$ws.Range("N716").Value = "$/malla 25 kilos"
$ws.Range("O716").Value = "Perú"
$ws.Range("P716").Value = 840
$ws.Range("Q716").Value = 25
$ws.Range("R716").Value = "Hortaliza"
